# Update Solar facility counts for 2023 and 2024 (column E) with
# updated data from upstream processes through 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 25 -> Open year 2023, column E -> Solar: 17 -> 19
$ws.Range("E25").Value = 19

# Row 26 -> Open year 2024, column E -> Solar: 12 -> 19
$ws.Range("E26").Value = 19
